# Lab 3 grading rubric edit
# - removes the stray "_GoBack" bookmark that used to sit on the
#   "Course number and Lab number" paragraph
# - gives the Reviewer row an explicit row height
# - rewords the Reviewer / Developer column headers
# - re-adds a "_GoBack" bookmark, now placed inside the word "form"
#   in the instructions paragraph ("for" + bookmark + "m for ")

$d = $word.ActiveDocument

# ---------------------------------------------------------------
# 1) Remove the "_GoBack" bookmark from the first table cell
#    ("Course number and Lab number"). Bookmark.Delete() is not
#    honoured by this host, so instead we rebuild the paragraph:
#    add a throw-away paragraph right after it (so the cell keeps a
#    trailing paragraph), delete the original (bookmark-carrying)
#    paragraph outright, then type the text back into what remains.
# ---------------------------------------------------------------
$table1 = $d.Tables.Item(1)

$firstCell = $table1.Cell(1, 1)
$firstPara = $firstCell.Range.Paragraphs.Item(1)
$tail = $firstPara.Range.Duplicate
$tail.Collapse(0)
$tail.InsertParagraphAfter()

$firstCellAgain = $table1.Cell(1, 1)
$firstCellAgain.Range.Paragraphs.Item(1).Range.Delete()

$firstCellFinal = $table1.Cell(1, 1)
$firstCellFinal.Range.Paragraphs.Item(1).Range.InsertBefore("Course number and Lab number")

# ---------------------------------------------------------------
# 2) Reviewer row gets an explicit height (305 twips = 15.25 pt)
# ---------------------------------------------------------------
$reviewerRow = $table1.Rows.Item(2)
$reviewerRow.Height = 15.25

# ---------------------------------------------------------------
# 3) Reword the Reviewer parenthetical
# ---------------------------------------------------------------
$reviewerRng = $d.Content
$reviewerRng.Find.Execute("You, the person doing the review)", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$reviewerRng.Text = "The person doing a review of the beta version)"

# ---------------------------------------------------------------
# 4) Reword the Developer parenthetical
# ---------------------------------------------------------------
$developerRng = $d.Content
$developerRng.Find.Execute("Developer (Person whose code you are reviewing)", $true, $false, $false, $false, `
    $false, $true, 1, $false, "", 0) | Out-Null
$developerRng.Text = "Developer (Person whose code is being reviewed)"

# ---------------------------------------------------------------
# 5) Re-insert the "_GoBack" bookmark inside "form for ", splitting
#    it into "for" + bookmark + "m for "
# ---------------------------------------------------------------
$formForRng = $d.Content
$formForRng.Find.Execute("form for ", $true, $false, $false, $false, $false, `
    $true, 1, $false, "", 0) | Out-Null
$bookmarkSpot = $d.Range($formForRng.Start + 3, $formForRng.Start + 3)
$d.Bookmarks.Add("_GoBack", $bookmarkSpot)
